# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets,
# mirroring the upstream logger's periodic dump.

$wb = $excel.ActiveWorkbook

$pirRows = @(
    @("2026-01-28", "16:20:26", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:20:26", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:20:31", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:20:36", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:20:41", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:20:46", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:20:51", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:20:56", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:21:01", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:21:06", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:21:11", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:21:16", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-28", "16:21:21", "16:00", "Bathroom", "No Motion", "Inactive")
)

$humidityRows = @(
    @("2026-01-28", "16:20:25", "16:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "16:20:28", "16:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "16:20:32", "16:00", "Bathroom", "88.4%", "Active"),
    @("2026-01-28", "16:20:40", "16:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "16:20:44", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:20:52", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:20:56", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:21:04", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:21:08", "16:00", "Bathroom", "87.4%", "Active"),
    @("2026-01-28", "16:21:12", "16:00", "Bathroom", "88.2%", "Active"),
    @("2026-01-28", "16:21:16", "16:00", "Bathroom", "88.3%", "Active"),
    @("2026-01-28", "16:21:20", "16:00", "Bathroom", "87.3%", "Active")
)

$temperatureRows = @(
    @("2026-01-28", "16:20:26", "16:00", "Bathroom", "22.7C", "Active"),
    @("2026-01-28", "16:20:29", "16:00", "Bathroom", "22.7C", "Active"),
    @("2026-01-28", "16:20:32", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:20:41", "16:00", "Bathroom", "22.7C", "Active"),
    @("2026-01-28", "16:20:45", "16:00", "Bathroom", "22.7C", "Active"),
    @("2026-01-28", "16:20:53", "16:00", "Bathroom", "22.7C", "Active"),
    @("2026-01-28", "16:20:57", "16:00", "Bathroom", "22.7C", "Active"),
    @("2026-01-28", "16:21:05", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:21:09", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:21:13", "16:00", "Bathroom", "22.7C", "Active"),
    @("2026-01-28", "16:21:17", "16:00", "Bathroom", "22.8C", "Active"),
    @("2026-01-28", "16:21:21", "16:00", "Bathroom", "22.7C", "Active")
)

function Append-LogRows($sheetName, $rows, $startRow, $percentCol) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = $startRow
    foreach ($row in $rows) {
        # Column A holds a plain-text date string ("2026-01-28"); force text
        # formatting first so Excel doesn't auto-coerce it into a date serial.
        $dateCell = $ws.Cells.Item($r, 1)
        $dateCell.NumberFormat = "@"
        $dateCell.Value = $row[0]
        $dateCell.Style = "Normal"

        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]

        if ($percentCol) {
            # Humidity values look like "88.2%"; force text so it isn't
            # auto-coerced into a numeric percentage.
            $valCell = $ws.Cells.Item($r, 5)
            $valCell.NumberFormat = "@"
            $valCell.Value = $row[4]
            $valCell.Style = "Normal"
        } else {
            $ws.Cells.Item($r, 5).Value = $row[4]
        }

        $ws.Cells.Item($r, 6).Value = $row[5]

        $r = $r + 1
    }
}

Append-LogRows "PIR" $pirRows 120 $false
Append-LogRows "Humidity" $humidityRows 119 $true
Append-LogRows "Temperature" $temperatureRows 119 $false
